# Updated tasks for Feb
$wb = $excel.ActiveWorkbook

# --- January sheet: selection moves from the Name column onto the Task column ---
$wsJan = $wb.Worksheets.Item("January")
$wsJan.Range("C3:C6").Select()

# --- February sheet: add the Task column (C) with the four task names ---
$wsFeb = $wb.Worksheets.Item("February")
$wsFeb.Activate()

$wsFeb.Range("C3").Value = "Migration Testing"
$wsFeb.Range("C4").Value = "Automation Testing"
$wsFeb.Range("C5").Value = "API Testing"
$wsFeb.Range("C6").Value = "Performance Testing"

# Match column width of the January sheet's Task column
$wsFeb.Columns.Item(3).ColumnWidth = 19.5703125

# Selection ends up parked just below the data, as in the source workbook
$wsFeb.Range("C10").Select()
